$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("B14").Value = 6352908
$ws.Range("E14").Value = 'Gangneung City'
$ws.Range("F14").Value = 'Siheung City AC'
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 'D'
$ws.Range("L14").Value = 2.5
$ws.Range("M14").Value = 3.2
$ws.Range("N14").Value = 2.5
$ws.Range("O14").Value = 2.625
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 2.5
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 1.95
$ws.Range("T14").Value = 1.85
$ws.Range("U14").Value = 2.25
$ws.Range("V14").Value = 2.025
$ws.Range("W14").Value = 1.775
$ws.Range("X14").Value = -1
$ws.Range("Y14").Value = 2
$ws.Range("Z14").Value = -1
$ws.Range("AA14").Value = 0
$ws.Range("AB14").Value = 0
$ws.Range("AC14").Value = -1
$ws.Range("AD14").Value = 0.7749999999999999

# Row 15
$ws.Range("B15").Value = 6352251
$ws.Range("E15").Value = 'Ulsan Citizen FC'
$ws.Range("F15").Value = 'Pocheon Citizen FC'
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 'A'
$ws.Range("L15").Value = 1.909
$ws.Range("M15").Value = 3.25
$ws.Range("N15").Value = 3.5
$ws.Range("O15").Value = 2.25
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 2.9
$ws.Range("R15").Value = -0.25
$ws.Range("S15").Value = 2.025
$ws.Range("T15").Value = 1.775
$ws.Range("U15").Value = 2
$ws.Range("V15").Value = 1.95
$ws.Range("W15").Value = 1.85
$ws.Range("X15").Value = -1
$ws.Range("Y15").Value = -1
$ws.Range("Z15").Value = 1.9
$ws.Range("AA15").Value = -1
$ws.Range("AB15").Value = 0.7749999999999999
$ws.Range("AC15").Value = -1
$ws.Range("AD15").Value = 0.8500000000000001

# Row 27
$ws.Range("B27").Value = 6352257
$ws.Range("E27").Value = 'Paju Citizen FC'
$ws.Range("F27").Value = 'Gimhae City'
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 2
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 'A'
$ws.Range("L27").Value = 2.25
$ws.Range("M27").Value = 3.1
$ws.Range("N27").Value = 2.875
$ws.Range("O27").Value = 2.15
$ws.Range("P27").Value = 3.2
$ws.Range("Q27").Value = 3.1
$ws.Range("R27").Value = -0.25
$ws.Range("S27").Value = 1.925
$ws.Range("T27").Value = 1.875
$ws.Range("U27").Value = 2
$ws.Range("V27").Value = 1.875
$ws.Range("W27").Value = 1.925
$ws.Range("X27").Value = -1
$ws.Range("Y27").Value = -1
$ws.Range("Z27").Value = 2.1
$ws.Range("AA27").Value = -1
$ws.Range("AB27").Value = 0.875
$ws.Range("AC27").Value = 0.875
$ws.Range("AD27").Value = -1

# Row 28
$ws.Range("B28").Value = 6352258
$ws.Range("E28").Value = 'Daejeon Korail'
$ws.Range("F28").Value = 'Siheung City AC'
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 1
$ws.Range("K28").Value = 'D'
$ws.Range("L28").Value = 2.375
$ws.Range("M28").Value = 3
$ws.Range("N28").Value = 2.75
$ws.Range("O28").Value = 2.375
$ws.Range("P28").Value = 3
$ws.Range("Q28").Value = 2.75
$ws.Range("R28").Value = -0.25
$ws.Range("S28").Value = 2.05
$ws.Range("T28").Value = 1.75
$ws.Range("U28").Value = 2.25
$ws.Range("V28").Value = 1.975
$ws.Range("W28").Value = 1.825
$ws.Range("X28").Value = -1
$ws.Range("Y28").Value = 2
$ws.Range("Z28").Value = -1
$ws.Range("AA28").Value = -0.5
$ws.Range("AB28").Value = 0.375
$ws.Range("AC28").Value = -0.5
$ws.Range("AD28").Value = 0.4125

# Row 35
$ws.Range("B35").Value = 6352922
$ws.Range("E35").Value = 'Mokpo City'
$ws.Range("F35").Value = 'Yangju Citizen'
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 'H'
$ws.Range("L35").Value = 1.666
$ws.Range("M35").Value = 3.4
$ws.Range("N35").Value = 4.5
$ws.Range("O35").Value = 1.666
$ws.Range("P35").Value = 3.5
$ws.Range("Q35").Value = 4.333
$ws.Range("R35").Value = -0.75
$ws.Range("S35").Value = 1.9
$ws.Range("T35").Value = 1.9
$ws.Range("U35").Value = 2.5
$ws.Range("V35").Value = 1.95
$ws.Range("W35").Value = 1.85
$ws.Range("X35").Value = 0.6659999999999999
$ws.Range("Y35").Value = -1
$ws.Range("Z35").Value = -1
$ws.Range("AA35").Value = 0.45
$ws.Range("AB35").Value = -0.5
$ws.Range("AC35").Value = -1
$ws.Range("AD35").Value = 0.8500000000000001

# Row 36
$ws.Range("B36").Value = 6352262
$ws.Range("E36").Value = 'Gangneung City'
$ws.Range("F36").Value = 'Pocheon Citizen FC'
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 3
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 2
$ws.Range("K36").Value = 'D'
$ws.Range("L36").Value = 2.4
$ws.Range("M36").Value = 3.2
$ws.Range("N36").Value = 2.625
$ws.Range("O36").Value = 2.375
$ws.Range("P36").Value = 3.2
$ws.Range("Q36").Value = 2.7
$ws.Range("R36").Value = 0
$ws.Range("S36").Value = 1.75
$ws.Range("T36").Value = 2.05
$ws.Range("U36").Value = 2
$ws.Range("V36").Value = 1.775
$ws.Range("W36").Value = 2.025
$ws.Range("X36").Value = -1
$ws.Range("Y36").Value = 2.2
$ws.Range("Z36").Value = -1
$ws.Range("AA36").Value = 0
$ws.Range("AB36").Value = 0
$ws.Range("AC36").Value = 0.7749999999999999
$ws.Range("AD36").Value = -1

# Row 49
$ws.Range("B49").Value = 6353327
$ws.Range("E49").Value = 'Yangpyeong FC'
$ws.Range("F49").Value = 'Gyeongju HNP'
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 1
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1
$ws.Range("K49").Value = 'A'
$ws.Range("L49").Value = 2.25
$ws.Range("M49").Value = 3.2
$ws.Range("N49").Value = 2.8
$ws.Range("O49").Value = 3.1
$ws.Range("P49").Value = 3.1
$ws.Range("Q49").Value = 2.1
$ws.Range("R49").Value = 0.25
$ws.Range("S49").Value = 1.875
$ws.Range("T49").Value = 1.925
$ws.Range("U49").Value = 2.25
$ws.Range("V49").Value = 1.975
$ws.Range("W49").Value = 1.725
$ws.Range("X49").Value = -1
$ws.Range("Y49").Value = -1
$ws.Range("Z49").Value = 1.1
$ws.Range("AA49").Value = -1
$ws.Range("AB49").Value = 0.925
$ws.Range("AC49").Value = -1
$ws.Range("AD49").Value = 0.7250000000000001

# Row 50
$ws.Range("B50").Value = 6352926
$ws.Range("E50").Value = 'Yangju Citizen'
$ws.Range("F50").Value = 'Siheung City AC'
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 1
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 1
$ws.Range("K50").Value = 'D'
$ws.Range("L50").Value = 2.35
$ws.Range("M50").Value = 3.25
$ws.Range("N50").Value = 2.625
$ws.Range("O50").Value = 2.875
$ws.Range("P50").Value = 3.3
$ws.Range("Q50").Value = 2.15
$ws.Range("R50").Value = 0.25
$ws.Range("S50").Value = 1.875
$ws.Range("T50").Value = 1.925
$ws.Range("U50").Value = 2.25
$ws.Range("V50").Value = 1.9
$ws.Range("W50").Value = 1.9
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = 2.3
$ws.Range("Z50").Value = -1
$ws.Range("AA50").Value = 0.4375
$ws.Range("AB50").Value = -0.5
$ws.Range("AC50").Value = -0.5
$ws.Range("AD50").Value = 0.45

# Row 72
$ws.Range("B72").Value = 6352279
$ws.Range("E72").Value = 'Yangju Citizen'
$ws.Range("F72").Value = 'Ulsan Citizen FC'
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 2
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 'A'
$ws.Range("L72").Value = 4.8
$ws.Range("M72").Value = 3.5
$ws.Range("N72").Value = 1.615
$ws.Range("O72").Value = 3.5
$ws.Range("P72").Value = 3.3
$ws.Range("Q72").Value = 1.95
$ws.Range("R72").Value = 0.5
$ws.Range("S72").Value = 1.825
$ws.Range("T72").Value = 1.975
$ws.Range("U72").Value = 2.25
$ws.Range("V72").Value = 1.95
$ws.Range("W72").Value = 1.85
$ws.Range("X72").Value = -1
$ws.Range("Y72").Value = -1
$ws.Range("Z72").Value = 0.95
$ws.Range("AA72").Value = -1
$ws.Range("AB72").Value = 0.9750000000000001
$ws.Range("AC72").Value = -0.5
$ws.Range("AD72").Value = 0.425

# Row 73
$ws.Range("B73").Value = 6353330
$ws.Range("E73").Value = 'Yangpyeong FC'
$ws.Range("F73").Value = 'Chuncheon FC'
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 'H'
$ws.Range("L73").Value = 4
$ws.Range("M73").Value = 3.1
$ws.Range("N73").Value = 1.833
$ws.Range("O73").Value = 2.375
$ws.Range("P73").Value = 2.875
$ws.Range("Q73").Value = 3.1
$ws.Range("R73").Value = 0
$ws.Range("S73").Value = 1.775
$ws.Range("T73").Value = 2.025
$ws.Range("U73").Value = 2
$ws.Range("V73").Value = 2.025
$ws.Range("W73").Value = 1.775
$ws.Range("X73").Value = 1.375
$ws.Range("Y73").Value = -1
$ws.Range("Z73").Value = -1
$ws.Range("AA73").Value = 0.7749999999999999
$ws.Range("AB73").Value = -1
$ws.Range("AC73").Value = -1
$ws.Range("AD73").Value = 0.7749999999999999

# Row 121
$ws.Range("B121").Value = 6352953
$ws.Range("E121").Value = 'Gimhae City'
$ws.Range("F121").Value = 'Daejeon Korail'
$ws.Range("G121").Value = 2
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 1
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 'H'
$ws.Range("L121").Value = 2.05
$ws.Range("M121").Value = 3.25
$ws.Range("N121").Value = 3.1
$ws.Range("O121").Value = 1.8
$ws.Range("P121").Value = 3.4
$ws.Range("Q121").Value = 3.75
$ws.Range("R121").Value = -0.5
$ws.Range("S121").Value = 1.825
$ws.Range("T121").Value = 1.975
$ws.Range("U121").Value = 2.5
$ws.Range("V121").Value = 2
$ws.Range("W121").Value = 1.8
$ws.Range("X121").Value = 0.8
$ws.Range("Y121").Value = -1
$ws.Range("Z121").Value = -1
$ws.Range("AA121").Value = 0.825
$ws.Range("AB121").Value = -1
$ws.Range("AC121").Value = -1
$ws.Range("AD121").Value = 0.8

# Row 122
$ws.Range("B122").Value = 6353334
$ws.Range("E122").Value = 'Gangneung City'
$ws.Range("F122").Value = 'Busan Trans Corp'
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 3
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2
$ws.Range("K122").Value = 'A'
$ws.Range("L122").Value = 2.1
$ws.Range("M122").Value = 3.3
$ws.Range("N122").Value = 3
$ws.Range("O122").Value = 2.375
$ws.Range("P122").Value = 3.3
$ws.Range("Q122").Value = 2.55
$ws.Range("R122").Value = 0
$ws.Range("S122").Value = 1.8
$ws.Range("T122").Value = 2
$ws.Range("U122").Value = 2.5
$ws.Range("V122").Value = 1.925
$ws.Range("W122").Value = 1.875
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = 1.55
$ws.Range("AA122").Value = -1
$ws.Range("AB122").Value = 1
$ws.Range("AC122").Value = 0.925
$ws.Range("AD122").Value = -1

# Row 123
$ws.Range("B123").Value = 6352956
$ws.Range("E123").Value = 'Ulsan Citizen FC'
$ws.Range("F123").Value = 'Yangpyeong FC'
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 'D'
$ws.Range("L123").Value = 2.05
$ws.Range("M123").Value = 3.5
$ws.Range("N123").Value = 2.9
$ws.Range("O123").Value = 2.75
$ws.Range("P123").Value = 3.5
$ws.Range("Q123").Value = 2.15
$ws.Range("R123").Value = 0.25
$ws.Range("S123").Value = 1.8
$ws.Range("T123").Value = 2
$ws.Range("U123").Value = 2.5
$ws.Range("V123").Value = 1.95
$ws.Range("W123").Value = 1.85
$ws.Range("X123").Value = -1
$ws.Range("Y123").Value = 2.5
$ws.Range("Z123").Value = -1
$ws.Range("AA123").Value = 0.4
$ws.Range("AB123").Value = -0.5
$ws.Range("AC123").Value = -1
$ws.Range("AD123").Value = 0.8500000000000001

# Row 124
$ws.Range("B124").Value = 6353335
$ws.Range("E124").Value = 'Siheung City AC'
$ws.Range("F124").Value = 'Chuncheon FC'
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 2
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 1
$ws.Range("K124").Value = 'D'
$ws.Range("L124").Value = 2.2
$ws.Range("M124").Value = 3.25
$ws.Range("N124").Value = 2.8
$ws.Range("O124").Value = 1.95
$ws.Range("P124").Value = 3.4
$ws.Range("Q124").Value = 3.2
$ws.Range("R124").Value = -0.5
$ws.Range("S124").Value = 2
$ws.Range("T124").Value = 1.8
$ws.Range("U124").Value = 2.25
$ws.Range("V124").Value = 1.775
$ws.Range("W124").Value = 2.025
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = 2.4
$ws.Range("Z124").Value = -1
$ws.Range("AA124").Value = -1
$ws.Range("AB124").Value = 0.8
$ws.Range("AC124").Value = 0.7749999999999999
$ws.Range("AD124").Value = -1

# Row 137
$ws.Range("B137").Value = 7867509
$ws.Range("E137").Value = 'Hwaseong FC'
$ws.Range("F137").Value = 'Gimhae City'
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1
$ws.Range("K137").Value = 'D'
$ws.Range("L137").Value = 2
$ws.Range("M137").Value = 3.25
$ws.Range("N137").Value = 3.25
$ws.Range("O137").Value = 2
$ws.Range("P137").Value = 3.1
$ws.Range("Q137").Value = 3.4
$ws.Range("R137").Value = -0.25
$ws.Range("S137").Value = 1.8
$ws.Range("T137").Value = 2
$ws.Range("U137").Value = 2.25
$ws.Range("V137").Value = 2.025
$ws.Range("W137").Value = 1.775
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 2.1
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = -0.5
$ws.Range("AB137").Value = 0.5
$ws.Range("AC137").Value = -0.5
$ws.Range("AD137").Value = 0.3875

# Row 138
$ws.Range("B138").Value = 7867506
$ws.Range("E138").Value = 'Ulsan Citizen FC'
$ws.Range("F138").Value = 'Yeoju'
$ws.Range("G138").Value = 2
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 1
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 'H'
$ws.Range("L138").Value = 1.8
$ws.Range("M138").Value = 3.25
$ws.Range("N138").Value = 4
$ws.Range("O138").Value = 1.8
$ws.Range("P138").Value = 3.25
$ws.Range("Q138").Value = 4
$ws.Range("R138").Value = -0.5
$ws.Range("S138").Value = 1.85
$ws.Range("T138").Value = 1.95
$ws.Range("U138").Value = 2.5
$ws.Range("V138").Value = 2
$ws.Range("W138").Value = 1.8
$ws.Range("X138").Value = 0.8
$ws.Range("Y138").Value = -1
$ws.Range("Z138").Value = -1
$ws.Range("AA138").Value = 0.8500000000000001
$ws.Range("AB138").Value = -1
$ws.Range("AC138").Value = -1
$ws.Range("AD138").Value = 0.8

# Row 139
$ws.Range("B139").Value = 7867508
$ws.Range("E139").Value = 'Chuncheon FC'
$ws.Range("F139").Value = 'Pocheon Citizen FC'
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = 3
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 'A'
$ws.Range("L139").Value = 2.1
$ws.Range("M139").Value = 3.25
$ws.Range("N139").Value = 3
$ws.Range("O139").Value = 2.1
$ws.Range("P139").Value = 3.25
$ws.Range("Q139").Value = 3
$ws.Range("R139").Value = -0.25
$ws.Range("S139").Value = 1.875
$ws.Range("T139").Value = 1.925
$ws.Range("U139").Value = 2.25
$ws.Range("V139").Value = 1.9
$ws.Range("W139").Value = 1.9
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = -1
$ws.Range("Z139").Value = 2
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = 0.925
$ws.Range("AC139").Value = 0.8999999999999999
$ws.Range("AD139").Value = -1

# Row 149
$ws.Range("B149").Value = 7867517
$ws.Range("E149").Value = 'Busan Trans Corp'
$ws.Range("F149").Value = 'Gangneung City'
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 3
$ws.Range("I149").Value = 0
$ws.Range("J149").Value = 1
$ws.Range("K149").Value = 'A'
$ws.Range("L149").Value = 2.375
$ws.Range("M149").Value = 3.1
$ws.Range("N149").Value = 2.7
$ws.Range("O149").Value = 2.6
$ws.Range("P149").Value = 3.2
$ws.Range("Q149").Value = 2.5
$ws.Range("R149").Value = 0
$ws.Range("S149").Value = 1.95
$ws.Range("T149").Value = 1.85
$ws.Range("U149").Value = 2.25
$ws.Range("V149").Value = 2.05
$ws.Range("W149").Value = 1.75
$ws.Range("X149").Value = -1
$ws.Range("Y149").Value = -1
$ws.Range("Z149").Value = 1.5
$ws.Range("AA149").Value = -1
$ws.Range("AB149").Value = 0.8500000000000001
$ws.Range("AC149").Value = 1.05
$ws.Range("AD149").Value = -1

# Row 150
$ws.Range("B150").Value = 7867516
$ws.Range("E150").Value = 'Yangpyeong FC'
$ws.Range("F150").Value = 'Chuncheon FC'
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 0
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 0
$ws.Range("K150").Value = 'H'
$ws.Range("L150").Value = 2.6
$ws.Range("M150").Value = 3.2
$ws.Range("N150").Value = 2.4
$ws.Range("O150").Value = 2.8
$ws.Range("P150").Value = 3.2
$ws.Range("Q150").Value = 2.25
$ws.Range("R150").Value = 0.25
$ws.Range("S150").Value = 1.775
$ws.Range("T150").Value = 2.025
$ws.Range("U150").Value = 2
$ws.Range("V150").Value = 1.85
$ws.Range("W150").Value = 1.95
$ws.Range("X150").Value = 1.8
$ws.Range("Y150").Value = -1
$ws.Range("Z150").Value = -1
$ws.Range("AA150").Value = 0.7749999999999999
$ws.Range("AB150").Value = -1
$ws.Range("AC150").Value = -1
$ws.Range("AD150").Value = 0.95

# Row 153
$ws.Range("B153").Value = 7867522
$ws.Range("E153").Value = 'Gangneung City'
$ws.Range("F153").Value = 'Pocheon Citizen FC'
$ws.Range("G153").Value = 3
$ws.Range("H153").Value = 1
$ws.Range("I153").Value = 2
$ws.Range("J153").Value = 1
$ws.Range("K153").Value = 'H'
$ws.Range("L153").Value = 2.2
$ws.Range("M153").Value = 3.2
$ws.Range("N153").Value = 2.9
$ws.Range("O153").Value = 1.95
$ws.Range("P153").Value = 3.3
$ws.Range("Q153").Value = 3.3
$ws.Range("R153").Value = -0.25
$ws.Range("S153").Value = 1.725
$ws.Range("T153").Value = 1.975
$ws.Range("U153").Value = 2
$ws.Range("V153").Value = 1.85
$ws.Range("W153").Value = 1.95
$ws.Range("X153").Value = 0.95
$ws.Range("Y153").Value = -1
$ws.Range("Z153").Value = -1
$ws.Range("AA153").Value = 0.7250000000000001
$ws.Range("AB153").Value = -1
$ws.Range("AC153").Value = 0.8500000000000001
$ws.Range("AD153").Value = -1

# Row 154
$ws.Range("B154").Value = 7873809
$ws.Range("E154").Value = 'Daegu FC Reserves'
$ws.Range("F154").Value = 'Ulsan Citizen FC'
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 2
$ws.Range("I154").Value = 0
$ws.Range("J154").Value = 1
$ws.Range("K154").Value = 'A'
$ws.Range("L154").Value = 2.75
$ws.Range("M154").Value = 4
$ws.Range("N154").Value = 2
$ws.Range("O154").Value = 2.6
$ws.Range("P154").Value = 3.5
$ws.Range("Q154").Value = 2.25
$ws.Range("R154").Value = 0.25
$ws.Range("S154").Value = 1.725
$ws.Range("T154").Value = 1.975
$ws.Range("U154").Value = 2.25
$ws.Range("V154").Value = 2
$ws.Range("W154").Value = 1.8
$ws.Range("X154").Value = -1
$ws.Range("Y154").Value = -1
$ws.Range("Z154").Value = 1.25
$ws.Range("AA154").Value = -1
$ws.Range("AB154").Value = 0.9750000000000001
$ws.Range("AC154").Value = 1
$ws.Range("AD154").Value = -1

# Row 167
$ws.Range("B167").Value = 7867537
$ws.Range("E167").Value = 'Hwaseong FC'
$ws.Range("F167").Value = 'Pocheon Citizen FC'
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 4
$ws.Range("I167").Value = 0
$ws.Range("J167").Value = 2
$ws.Range("K167").Value = 'A'
$ws.Range("L167").Value = 1.571
$ws.Range("M167").Value = 3.6
$ws.Range("N167").Value = 5
$ws.Range("O167").Value = 1.533
$ws.Range("P167").Value = 3.6
$ws.Range("Q167").Value = 5.25
$ws.Range("R167").Value = -1
$ws.Range("S167").Value = 1.9
$ws.Range("T167").Value = 1.9
$ws.Range("U167").Value = 2.5
$ws.Range("V167").Value = 1.875
$ws.Range("W167").Value = 1.925
$ws.Range("X167").Value = -1
$ws.Range("Y167").Value = -1
$ws.Range("Z167").Value = 4.25
$ws.Range("AA167").Value = -1
$ws.Range("AB167").Value = 0.8999999999999999
$ws.Range("AC167").Value = 0.875
$ws.Range("AD167").Value = -1

# Row 168
$ws.Range("B168").Value = 7867538
$ws.Range("E168").Value = 'Gangneung City'
$ws.Range("F168").Value = 'Yeoju'
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0
$ws.Range("I168").Value = 0
$ws.Range("J168").Value = 0
$ws.Range("K168").Value = 'D'
$ws.Range("L168").Value = 1.4
$ws.Range("M168").Value = 4.333
$ws.Range("N168").Value = 6
$ws.Range("O168").Value = 1.5
$ws.Range("P168").Value = 4
$ws.Range("Q168").Value = 5.25
$ws.Range("R168").Value = -1
$ws.Range("S168").Value = 1.9
$ws.Range("T168").Value = 1.9
$ws.Range("U168").Value = 2.5
$ws.Range("V168").Value = 1.975
$ws.Range("W168").Value = 1.825
$ws.Range("X168").Value = -1
$ws.Range("Y168").Value = 3
$ws.Range("Z168").Value = -1
$ws.Range("AA168").Value = -1
$ws.Range("AB168").Value = 0.8999999999999999
$ws.Range("AC168").Value = -1
$ws.Range("AD168").Value = 0.825

# Row 169
$ws.Range("B169").Value = 7867539
$ws.Range("E169").Value = 'Gyeongju HNP'
$ws.Range("F169").Value = 'Mokpo City'
$ws.Range("G169").Value = 2
$ws.Range("H169").Value = 0
$ws.Range("I169").Value = 2
$ws.Range("J169").Value = 0
$ws.Range("K169").Value = 'H'
$ws.Range("L169").Value = 2
$ws.Range("M169").Value = 3.6
$ws.Range("N169").Value = 3
$ws.Range("O169").Value = 1.444
$ws.Range("P169").Value = 4.333
$ws.Range("Q169").Value = 5.75
$ws.Range("R169").Value = -1.25
$ws.Range("S169").Value = 2.025
$ws.Range("T169").Value = 1.775
$ws.Range("U169").Value = 2.5
$ws.Range("V169").Value = 1.9
$ws.Range("W169").Value = 1.9
$ws.Range("X169").Value = 0.444
$ws.Range("Y169").Value = -1
$ws.Range("Z169").Value = -1
$ws.Range("AA169").Value = 1.025
$ws.Range("AB169").Value = -1
$ws.Range("AC169").Value = -1
$ws.Range("AD169").Value = 0.8999999999999999

# Row 170
$ws.Range("B170").Value = 7867540
$ws.Range("E170").Value = 'Changwon City'
$ws.Range("F170").Value = 'Yangpyeong FC'
$ws.Range("G170").Value = 2
$ws.Range("H170").Value = 0
$ws.Range("I170").Value = 2
$ws.Range("J170").Value = 0
$ws.Range("K170").Value = 'H'
$ws.Range("L170").Value = 2
$ws.Range("M170").Value = 3.3
$ws.Range("N170").Value = 3.25
$ws.Range("O170").Value = 2
$ws.Range("P170").Value = 3.3
$ws.Range("Q170").Value = 3.25
$ws.Range("R170").Value = -0.25
$ws.Range("S170").Value = 1.8
$ws.Range("T170").Value = 2
$ws.Range("U170").Value = 2.25
$ws.Range("V170").Value = 2
$ws.Range("W170").Value = 1.8
$ws.Range("X170").Value = 1
$ws.Range("Y170").Value = -1
$ws.Range("Z170").Value = -1
$ws.Range("AA170").Value = 0.8
$ws.Range("AB170").Value = -1
$ws.Range("AC170").Value = -0.5
$ws.Range("AD170").Value = 0.4

# Row 175
$ws.Range("B175").Value = 7867546
$ws.Range("E175").Value = 'Siheung City AC'
$ws.Range("F175").Value = 'Chuncheon FC'
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0
$ws.Range("I175").Value = 0
$ws.Range("J175").Value = 0
$ws.Range("K175").Value = 'D'
$ws.Range("L175").Value = 1.8
$ws.Range("M175").Value = 3.2
$ws.Range("N175").Value = 4
$ws.Range("O175").Value = 1.8
$ws.Range("P175").Value = 3.2
$ws.Range("Q175").Value = 4
$ws.Range("R175").Value = -0.5
$ws.Range("S175").Value = 1.85
$ws.Range("T175").Value = 1.95
$ws.Range("U175").Value = 2.25
$ws.Range("V175").Value = 1.8
$ws.Range("W175").Value = 2
$ws.Range("X175").Value = -1
$ws.Range("Y175").Value = 2.2
$ws.Range("Z175").Value = -1
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = 0.95
$ws.Range("AC175").Value = -1
$ws.Range("AD175").Value = 1

# Row 176
$ws.Range("B176").Value = 7867544
$ws.Range("E176").Value = 'Yangpyeong FC'
$ws.Range("F176").Value = 'Hwaseong FC'
$ws.Range("G176").Value = 2
$ws.Range("H176").Value = 1
$ws.Range("I176").Value = 1
$ws.Range("J176").Value = 0
$ws.Range("K176").Value = 'H'
$ws.Range("L176").Value = 4
$ws.Range("M176").Value = 3.75
$ws.Range("N176").Value = 1.666
$ws.Range("O176").Value = 5.5
$ws.Range("P176").Value = 3.4
$ws.Range("Q176").Value = 1.615
$ws.Range("R176").Value = 0.75
$ws.Range("S176").Value = 2
$ws.Range("T176").Value = 1.8
$ws.Range("U176").Value = 2.25
$ws.Range("V176").Value = 1.975
$ws.Range("W176").Value = 1.825
$ws.Range("X176").Value = 4.5
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = 1
$ws.Range("AB176").Value = -1
$ws.Range("AC176").Value = 0.9750000000000001
$ws.Range("AD176").Value = -1

# Row 177
$ws.Range("B177").Value = 7867545
$ws.Range("E177").Value = 'Gyeongju HNP'
$ws.Range("F177").Value = 'Ulsan Citizen FC'
$ws.Range("G177").Value = 1
$ws.Range("H177").Value = 2
$ws.Range("I177").Value = 1
$ws.Range("J177").Value = 1
$ws.Range("K177").Value = 'A'
$ws.Range("L177").Value = 2.15
$ws.Range("M177").Value = 3.1
$ws.Range("N177").Value = 3
$ws.Range("O177").Value = 1.727
$ws.Range("P177").Value = 3.4
$ws.Range("Q177").Value = 4.5
$ws.Range("R177").Value = -0.75
$ws.Range("S177").Value = 1.95
$ws.Range("T177").Value = 1.85
$ws.Range("U177").Value = 2.25
$ws.Range("V177").Value = 1.9
$ws.Range("W177").Value = 1.9
$ws.Range("X177").Value = -1
$ws.Range("Y177").Value = -1
$ws.Range("Z177").Value = 3.5
$ws.Range("AA177").Value = -1
$ws.Range("AB177").Value = 0.8500000000000001
$ws.Range("AC177").Value = 0.8999999999999999
$ws.Range("AD177").Value = -1

# Row 187
$ws.Range("B187").Value = 7867557
$ws.Range("E187").Value = 'Chuncheon FC'
$ws.Range("F187").Value = 'Gyeongju HNP'
$ws.Range("G187").Value = 3
$ws.Range("H187").Value = 2
$ws.Range("I187").Value = 2
$ws.Range("J187").Value = 0
$ws.Range("K187").Value = 'H'
$ws.Range("L187").Value = 2.5
$ws.Range("M187").Value = 2.875
$ws.Range("N187").Value = 2.75
$ws.Range("O187").Value = 3.8
$ws.Range("P187").Value = 3
$ws.Range("Q187").Value = 1.95
$ws.Range("R187").Value = 0.5
$ws.Range("S187").Value = 1.825
$ws.Range("T187").Value = 1.975
$ws.Range("U187").Value = 2
$ws.Range("V187").Value = 1.775
$ws.Range("W187").Value = 2.025
$ws.Range("X187").Value = 2.8
$ws.Range("Y187").Value = -1
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = 0.825
$ws.Range("AB187").Value = -1
$ws.Range("AC187").Value = 0.7749999999999999
$ws.Range("AD187").Value = -1

# Row 188
$ws.Range("B188").Value = 7867555
$ws.Range("E188").Value = 'Siheung City AC'
$ws.Range("F188").Value = 'Gangneung City'
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 2
$ws.Range("I188").Value = 0
$ws.Range("J188").Value = 0
$ws.Range("K188").Value = 'A'
$ws.Range("L188").Value = 1.909
$ws.Range("M188").Value = 2.9
$ws.Range("N188").Value = 4
$ws.Range("O188").Value = 2.15
$ws.Range("P188").Value = 2.8
$ws.Range("Q188").Value = 3.3
$ws.Range("R188").Value = -0.25
$ws.Range("S188").Value = 1.85
$ws.Range("T188").Value = 1.95
$ws.Range("U188").Value = 2.25
$ws.Range("V188").Value = 1.975
$ws.Range("W188").Value = 1.825
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = -1
$ws.Range("Z188").Value = 2.3
$ws.Range("AA188").Value = -1
$ws.Range("AB188").Value = 0.95
$ws.Range("AC188").Value = -0.5
$ws.Range("AD188").Value = 0.4125

# Row 199
$ws.Range("B199").Value = 7873803
$ws.Range("E199").Value = 'Gyeongju HNP'
$ws.Range("F199").Value = 'Daegu FC Reserves'
$ws.Range("G199").Value = 3
$ws.Range("H199").Value = 0
$ws.Range("I199").Value = 1
$ws.Range("J199").Value = 0
$ws.Range("K199").Value = 'H'
$ws.Range("L199").Value = 1.333
$ws.Range("M199").Value = 4.5
$ws.Range("N199").Value = 7
$ws.Range("O199").Value = 1.4
$ws.Range("P199").Value = 4.333
$ws.Range("Q199").Value = 5.75
$ws.Range("R199").Value = -1.25
$ws.Range("S199").Value = 1.85
$ws.Range("T199").Value = 1.95
$ws.Range("U199").Value = 2.75
$ws.Range("V199").Value = 1.75
$ws.Range("W199").Value = 1.95
$ws.Range("X199").Value = 0.3999999999999999
$ws.Range("Y199").Value = -1
$ws.Range("Z199").Value = -1
$ws.Range("AA199").Value = 0.8500000000000001
$ws.Range("AB199").Value = -1
$ws.Range("AC199").Value = 0.375
$ws.Range("AD199").Value = -0.5

# Row 201
$ws.Range("B201").Value = 7867567
$ws.Range("E201").Value = 'Gimhae City'
$ws.Range("F201").Value = 'Mokpo City'
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 2
$ws.Range("I201").Value = 0
$ws.Range("J201").Value = 0
$ws.Range("K201").Value = 'A'
$ws.Range("L201").Value = 1.571
$ws.Range("M201").Value = 3.6
$ws.Range("N201").Value = 5
$ws.Range("O201").Value = 1.65
$ws.Range("P201").Value = 3.5
$ws.Range("Q201").Value = 4.5
$ws.Range("R201").Value = -0.75
$ws.Range("S201").Value = 1.9
$ws.Range("T201").Value = 1.9
$ws.Range("U201").Value = 2.25
$ws.Range("V201").Value = 1.825
$ws.Range("W201").Value = 1.975
$ws.Range("X201").Value = -1
$ws.Range("Y201").Value = -1
$ws.Range("Z201").Value = 3.5
$ws.Range("AA201").Value = -1
$ws.Range("AB201").Value = 0.8999999999999999
$ws.Range("AC201").Value = -0.5
$ws.Range("AD201").Value = 0.4875
